$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All existing data rows (2-17) share identical formatting, so rather than
# using Rows.Insert (which would inherit the bold header row's style for
# the newly-inserted row), shift the existing data down one row by copying
# values from the bottom up (cell by cell, via Value2), then write the new
# data point into row 2. This leaves the existing per-row formatting in
# rows 2-17 untouched, and the row that is freshly populated at the bottom
# (row 18) has its date-time columns' number format re-applied to match.
for ($r = 17; $r -ge 2; $r--) {
    for ($c = 1; $c -le 6; $c++) {
        $v = $ws.Cells.Item($r, $c).Value2
        $ws.Cells.Item($r + 1, $c).Value2 = $v
    }
}

# Row 18 is brand-new territory (previously unused), so give its date
# columns (A:C) the same number format used by the rest of the table.
$ws.Range("A18:C18").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row 2 with the latest ISS Rt data point (Dec 4th update)
$ws.Cells.Item(2, 1).Value2 = 44146
$ws.Cells.Item(2, 2).Value2 = 44159
$ws.Cells.Item(2, 3).Value2 = 44152.99998842592
$ws.Cells.Item(2, 4).Value2 = 0.91
$ws.Cells.Item(2, 5).Value2 = 0.79
$ws.Cells.Item(2, 6).Value2 = 1.08
